$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Cells.Item(80, 8).Value = 2044.4
$ws.Cells.Item(80, 9).Value = 2274
$ws.Cells.Item(80, 10).Value = 1700
$ws.Cells.Item(80, 11).Value = 6822
$ws.Cells.Item(80, 12).Value = 5100
$ws.Cells.Item(80, 13).Value = -5824
$ws.Cells.Item(80, 14).Value = -7096
# Row 83
$ws.Cells.Item(83, 8).Value = 2044.4
$ws.Cells.Item(83, 9).Value = 2274
$ws.Cells.Item(83, 10).Value = 1700
$ws.Cells.Item(83, 11).Value = 20466
$ws.Cells.Item(83, 12).Value = 15300
$ws.Cells.Item(83, 13).Value = -15474
$ws.Cells.Item(83, 14).Value = -25284
# Row 98
$ws.Cells.Item(98, 8).Value = 1744.6364
$ws.Cells.Item(98, 9).Value = 1744.1
$ws.Cells.Item(98, 11).Value = 1744.1
$ws.Cells.Item(98, 13).Value = -246.0999999999999
# Row 100
$ws.Cells.Item(100, 8).Value = 1923.625
$ws.Cells.Item(100, 9).Value = 1847.25
$ws.Cells.Item(100, 11).Value = 1847.25
$ws.Cells.Item(100, 13).Value = -1306.25
# Row 122
$ws.Cells.Item(122, 8).Value = 1744.6364
$ws.Cells.Item(122, 9).Value = 1744.1
$ws.Cells.Item(122, 11).Value = 5232.299999999999
$ws.Cells.Item(122, 13).Value = -2782.299999999999
# Row 129
$ws.Cells.Item(129, 8).Value = 2377
$ws.Cells.Item(129, 9).Value = 839.6667
$ws.Cells.Item(129, 10).Value = 2953.5
$ws.Cells.Item(129, 11).Value = 2519.0001
$ws.Cells.Item(129, 12).Value = 8860.5
$ws.Cells.Item(129, 13).Value = 2480.9999
$ws.Cells.Item(129, 14).Value = -18860.5
# Row 135
$ws.Cells.Item(135, 8).Value = 2001.0667
$ws.Cells.Item(135, 9).Value = 1813.2222
$ws.Cells.Item(135, 11).Value = 16318.9998
$ws.Cells.Item(135, 13).Value = -13783.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6619.5835
$ws.Cells.Item(32, 9).Value = 4939.8945
$ws.Cells.Item(32, 10).Value = 13002.4
$ws.Cells.Item(32, 11).Value = 4939.8945
$ws.Cells.Item(32, 12).Value = 13002.4
$ws.Cells.Item(32, 13).Value = -4652.8945
$ws.Cells.Item(32, 14).Value = -13576.4
# Row 61
$ws.Cells.Item(61, 8).Value = 3695.7368
$ws.Cells.Item(61, 9).Value = 3723.2778
$ws.Cells.Item(61, 10).Value = 3200
$ws.Cells.Item(61, 11).Value = 3723.2778
$ws.Cells.Item(61, 12).Value = 3200
$ws.Cells.Item(61, 13).Value = -3511.2778
$ws.Cells.Item(61, 14).Value = -3624
# Row 122
$ws.Cells.Item(122, 8).Value = 486327.9
$ws.Cells.Item(122, 9).Value = 630566.7
$ws.Cells.Item(122, 11).Value = 1891700.1
$ws.Cells.Item(122, 13).Value = -1889250.1
# Row 132
$ws.Cells.Item(132, 8).Value = 4366.1333
$ws.Cells.Item(132, 9).Value = 5665.3335
$ws.Cells.Item(132, 11).Value = 16996.0005
$ws.Cells.Item(132, 13).Value = -14466.0005
# Row 133
$ws.Cells.Item(133, 8).Value = 70000
$ws.Cells.Item(133, 10).Value = 70000
$ws.Cells.Item(133, 12).Value = 70000
$ws.Cells.Item(133, 14).Value = -75060
# Row 136
$ws.Cells.Item(136, 8).Value = 3695.7368
$ws.Cells.Item(136, 9).Value = 3723.2778
$ws.Cells.Item(136, 10).Value = 3200
$ws.Cells.Item(136, 11).Value = 11169.8334
$ws.Cells.Item(136, 12).Value = 9600
$ws.Cells.Item(136, 13).Value = -8619.8334
$ws.Cells.Item(136, 14).Value = -14700

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Cells.Item(12, 8).Value = 1143.5454
$ws.Cells.Item(12, 9).Value = 471.5
$ws.Cells.Item(12, 10).Value = 1527.5714
$ws.Cells.Item(12, 11).Value = 471.5
$ws.Cells.Item(12, 12).Value = 1527.5714
$ws.Cells.Item(12, 13).Value = -303.5
$ws.Cells.Item(12, 14).Value = -1863.5714
# Row 86
$ws.Cells.Item(86, 8).Value = 3034.2856
$ws.Cells.Item(86, 9).Value = 3286.2
$ws.Cells.Item(86, 11).Value = 3286.2
$ws.Cells.Item(86, 13).Value = -2163.2
# Row 89
$ws.Cells.Item(89, 8).Value = 3034.2856
$ws.Cells.Item(89, 9).Value = 3286.2
$ws.Cells.Item(89, 11).Value = 16431
$ws.Cells.Item(89, 13).Value = -10815

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Cells.Item(107, 8).Value = 12500971
$ws.Cells.Item(107, 9).Value = 19231358
$ws.Cells.Item(107, 10).Value = 1678.5714
$ws.Cells.Item(107, 11).Value = 19231358
$ws.Cells.Item(107, 12).Value = 1678.5714
$ws.Cells.Item(107, 13).Value = -19229438
$ws.Cells.Item(107, 14).Value = -5518.5714

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Cells.Item(4, 8).Value = 33103418
$ws.Cells.Item(4, 10).Value = 7000
$ws.Cells.Item(4, 12).Value = 21000
$ws.Cells.Item(4, 14).Value = -21224
# Row 68
$ws.Cells.Item(68, 8).Value = 901
$ws.Cells.Item(68, 9).Value = 850.75
$ws.Cells.Item(68, 10).Value = 1001.5
$ws.Cells.Item(68, 11).Value = 2552.25
$ws.Cells.Item(68, 12).Value = 3004.5
$ws.Cells.Item(68, 13).Value = -1741.25
$ws.Cells.Item(68, 14).Value = -4626.5
# Row 71
$ws.Cells.Item(71, 8).Value = 901
$ws.Cells.Item(71, 9).Value = 850.75
$ws.Cells.Item(71, 10).Value = 1001.5
$ws.Cells.Item(71, 11).Value = 7656.75
$ws.Cells.Item(71, 12).Value = 9013.5
$ws.Cells.Item(71, 13).Value = -3600.75
$ws.Cells.Item(71, 14).Value = -17125.5
# Row 132
$ws.Cells.Item(132, 8).Value = 4518.8
$ws.Cells.Item(132, 9).Value = 4497
$ws.Cells.Item(132, 10).Value = 4533.3335
$ws.Cells.Item(132, 11).Value = 40473
$ws.Cells.Item(132, 12).Value = 40800.0015
$ws.Cells.Item(132, 13).Value = -37943
$ws.Cells.Item(132, 14).Value = -45860.0015

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 7483.1665
$ws.Cells.Item(70, 9).Value = 7224.75
$ws.Cells.Item(70, 11).Value = 7224.75
$ws.Cells.Item(70, 13).Value = -6954.75
# Row 73
$ws.Cells.Item(73, 8).Value = 7483.1665
$ws.Cells.Item(73, 9).Value = 7224.75
$ws.Cells.Item(73, 11).Value = 7224.75
$ws.Cells.Item(73, 13).Value = -6288.75
# Row 80
$ws.Cells.Item(80, 8).Value = 4750
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
# Row 83
$ws.Cells.Item(83, 8).Value = 4750
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
# Row 97
$ws.Cells.Item(97, 8).Value = 2147.5715
$ws.Cells.Item(97, 9).Value = 1539
$ws.Cells.Item(97, 10).Value = 3669
$ws.Cells.Item(97, 11).Value = 1539
$ws.Cells.Item(97, 12).Value = 3669
$ws.Cells.Item(97, 13).Value = -1043
$ws.Cells.Item(97, 14).Value = -4661
# Row 102
$ws.Cells.Item(102, 8).Value = 2917.8
$ws.Cells.Item(102, 9).Value = 2600
$ws.Cells.Item(102, 11).Value = 2600
$ws.Cells.Item(102, 13).Value = -978
# Row 122
$ws.Cells.Item(122, 8).Value = 61359.53
$ws.Cells.Item(122, 9).Value = 2057.6365
$ws.Cells.Item(122, 10).Value = 170079.67
$ws.Cells.Item(122, 11).Value = 6172.9095
$ws.Cells.Item(122, 12).Value = 510239.01
$ws.Cells.Item(122, 13).Value = -3722.9095
$ws.Cells.Item(122, 14).Value = -515139.01
# Row 126
$ws.Cells.Item(126, 8).Value = 1102.25
$ws.Cells.Item(126, 9).Value = 1037
$ws.Cells.Item(126, 11).Value = 3111
$ws.Cells.Item(126, 13).Value = -641
# Row 132
$ws.Cells.Item(132, 8).Value = 1713.5714
$ws.Cells.Item(132, 9).Value = 1713.5714
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 5140.7142
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2610.7142
$ws.Cells.Item(132, 14).ClearContents()
# Row 135
$ws.Cells.Item(135, 8).Value = 101229.7
$ws.Cells.Item(135, 9).Value = 101060.336
$ws.Cells.Item(135, 10).Value = 102500
$ws.Cells.Item(135, 11).Value = 101060.336
$ws.Cells.Item(135, 12).Value = 102500
$ws.Cells.Item(135, 13).Value = -95990.336
$ws.Cells.Item(135, 14).Value = -112640

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 5136.75
$ws.Cells.Item(7, 9).Value = 5136.75
$ws.Cells.Item(7, 11).Value = 5136.75
$ws.Cells.Item(7, 13).Value = -5024.75
# Row 46
$ws.Cells.Item(46, 8).Value = 2284.238
$ws.Cells.Item(46, 10).Value = 2286.2666
$ws.Cells.Item(46, 12).Value = 2286.2666
$ws.Cells.Item(46, 14).Value = -2662.2666
# Row 122
$ws.Cells.Item(122, 8).Value = 1383.3334
$ws.Cells.Item(122, 9).Value = 1383.3334
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4150.0002
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1700.0002
$ws.Cells.Item(122, 14).ClearContents()
# Row 126
$ws.Cells.Item(126, 8).Value = 5136.75
$ws.Cells.Item(126, 9).Value = 5136.75
$ws.Cells.Item(126, 11).Value = 15410.25
$ws.Cells.Item(126, 13).Value = -12940.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Cells.Item(4, 8).Value = 1543788.5
$ws.Cells.Item(4, 10).Value = 1750
$ws.Cells.Item(4, 12).Value = 1750
$ws.Cells.Item(4, 14).Value = -1976
# Row 126
$ws.Cells.Item(126, 8).Value = 2072.25
$ws.Cells.Item(126, 9).Value = 1797
$ws.Cells.Item(126, 10).Value = 2898
$ws.Cells.Item(126, 11).Value = 5391
$ws.Cells.Item(126, 12).Value = 8694
$ws.Cells.Item(126, 13).Value = -2921
$ws.Cells.Item(126, 14).Value = -13634
# Row 132
$ws.Cells.Item(132, 8).Value = 2987.7144
$ws.Cells.Item(132, 9).Value = 2643.353
$ws.Cells.Item(132, 11).Value = 7930.059
$ws.Cells.Item(132, 13).Value = -5400.059

Write-Output "Applied all Seraph_Profits market-data updates."